# Add a new "Appointment notes" column header to the patient list sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H (skipping G) gets a new header cell, matching the target edit.
$ws.Range("H1").Value = "Appointment notes"
